$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.747.60"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.724.67"
$ws.Range("D4").Value = "0.9976"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "241.10"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "0.4845"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "0.06196"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "1.728.19"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "15.98"
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").Value = "0.06901"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "0.6077"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "0.9982"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "26.562.45"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "0.9976"
$ws.Range("D19").Value = "0.000007149"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "1.950.07"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "4.434"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "8.556"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "5.069"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "136.97"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "15.26"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "1.770"
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "105.90"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Value = "3.937"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "0.07947"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "3.695"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "0.04492"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.597"
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.009"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6200"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "0.9277"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.032"
$ws.Range("E38").Value = "  +4.00%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.436"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "0.9973"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01496"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.660"
$ws.Range("E42").Value = "  +5.97%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "99.57"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3837"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "6.856"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1156"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05392"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "7.891"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.11"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "51.51"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.228"
$ws.Range("E51").Value = "  +0.00%  "
